$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to text before writing, so numeric-looking
# strings like "1.00" or "380.10" are preserved verbatim instead of being
# auto-coerced to numbers by Excel's type inference.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "64.767.61"
$ws.Range("E2").Value = "  -1.92%  "

$ws.Range("D3").Value = "3.428.13"
$ws.Range("E3").Value = "  -2.84%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "580.82"
$ws.Range("E5").Value = "  -3.81%  "

$ws.Range("D6").Value = "133.78"
$ws.Range("E6").Value = "  -7.05%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "3.424.03"
$ws.Range("E8").Value = "  -2.98%  "

$ws.Range("E9").Value = "  -5.62%  "

$ws.Range("E10").Value = "  -8.15%  "

$ws.Range("D11").Value = "7.09"
$ws.Range("E11").Value = "  -9.30%  "

$ws.Range("D12").Value = "0.373"
$ws.Range("E12").Value = "  -8.30%  "

$ws.Range("D13").Value = "4.012.29"
$ws.Range("E13").Value = "  -2.87%  "

$ws.Range("D14").Value = "0.0000178"
$ws.Range("E14").Value = "  -7.98%  "

$ws.Range("D15").Value = "3.439.94"
$ws.Range("E15").Value = "  -2.53%  "

$ws.Range("E16").Value = "  -1.70%  "

$ws.Range("D17").Value = "26.12"
$ws.Range("E17").Value = "  -7.90%  "

$ws.Range("D18").Value = "64.760.31"
$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("D19").Value = "9.46"
$ws.Range("E19").Value = "  -14.06%  "

$ws.Range("D20").Value = "5.73"
$ws.Range("E20").Value = "  -7.35%  "

$ws.Range("D21").Value = "13.45"
$ws.Range("E21").Value = "  -7.67%  "

$ws.Range("D22").Value = "380.10"
$ws.Range("E22").Value = "  -9.67%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "0.540"
$ws.Range("E24").Value = "  -9.11%  "

$ws.Range("D25").Value = "71.77"
$ws.Range("E25").Value = "  -6.72%  "

$ws.Range("D26").Value = "3.566.20"
$ws.Range("E26").Value = "  -2.89%  "

$ws.Range("D27").Value = "0.0000104"

$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("D29").Value = "7.13"
$ws.Range("E29").Value = "  -8.87%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.18"
$ws.Range("E30").Value = "  -11.60%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "7.94"
$ws.Range("E31").Value = "  -10.64%  "

$ws.Range("D32").Value = "3.447.11"
$ws.Range("E32").Value = "  -2.56%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "0.142"
$ws.Range("E34").Value = "  -8.91%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "22.87"
$ws.Range("E35").Value = "  -5.34%  "

$ws.Range("D36").Value = "169.50"
$ws.Range("E36").Value = "  -4.88%  "

$ws.Range("D37").Value = "1.17"
$ws.Range("E37").Value = "  -12.91%  "

$ws.Range("D38").Value = "6.67"
$ws.Range("E38").Value = "  -11.86%  "

$ws.Range("D39").Value = "1.44"
$ws.Range("E39").Value = "  -11.67%  "

$ws.Range("D40").Value = "4.54"
$ws.Range("E40").Value = "  -13.18%  "

$ws.Range("D41").Value = "0.0756"
$ws.Range("E41").Value = "  -7.53%  "

$ws.Range("D42").Value = "0.801"
$ws.Range("E42").Value = "  -6.62%  "

$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "41.90"
$ws.Range("E44").Value = "  -7.55%  "

$ws.Range("D45").Value = "4.27"
$ws.Range("E45").Value = "  -14.18%  "

$ws.Range("D46").Value = "1.61"
$ws.Range("E46").Value = "  -8.94%  "

$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").Value = "22.39"
$ws.Range("E48").Value = "  -5.59%  "

$ws.Range("E49").Value = "  -8.43%  "

$ws.Range("D50").Value = "2.196.31"
$ws.Range("E50").Value = "  -5.53%  "

$ws.Range("D51").Value = "1.97"
$ws.Range("E51").Value = "  -17.75%  "

# Restore the original "General" formatting/style now that the literal text
# has been committed to each cell.
$dataRange.NumberFormat = "General"
$dataRange.Style = "Normal"
